$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "319.70"
Set-TextValue "E2" "5.52%"
Set-TextValue "G2" "13"
Set-TextValue "D3" "48.81"
Set-TextValue "E3" "13.98%"
Set-TextValue "G3" "13"
Set-TextValue "D4" "5.282"
Set-TextValue "E4" "4.84%"
Set-TextValue "G4" "13"
Set-TextValue "D5" "0.08071"
Set-TextValue "E5" "5.26%"
Set-TextValue "G5" "13"
Set-TextValue "D6" "4.571"
Set-TextValue "E6" "3.75%"
Set-TextValue "G6" "13"
Set-TextValue "D7" "1.297"
Set-TextValue "E7" "27.31%"
Set-TextValue "G7" "13"
Set-TextValue "D8" "1.649"
Set-TextValue "E8" "2.09%"
Set-TextValue "G8" "13"
Set-TextValue "D9" "0.1298"
Set-TextValue "E9" "6.31%"
Set-TextValue "G9" "13"
Set-TextValue "E10" "3.88%"
Set-TextValue "G10" "13"
Set-TextValue "D11" "0.09413"
Set-TextValue "E11" "2.35%"
Set-TextValue "G11" "13"
Set-TextValue "D12" "0.04596"
Set-TextValue "E12" "10.64%"
Set-TextValue "G12" "13"
Set-TextValue "D13" "0.1045"
Set-TextValue "E13" "0.14%"
Set-TextValue "G13" "13"
Set-TextValue "D14" "0.001331"
Set-TextValue "E14" "3.94%"
Set-TextValue "G14" "13"
Set-TextValue "D15" "0.04167"
Set-TextValue "E15" "0.73%"
Set-TextValue "G15" "13"
Set-TextValue "D16" "0.005830"
Set-TextValue "E16" "-1.98%"
Set-TextValue "G16" "13"
Set-TextValue "D17" "3.338"
Set-TextValue "E17" "0.52%"
Set-TextValue "G17" "13"
Set-TextValue "D18" "2.438"
Set-TextValue "E18" "2.27%"
Set-TextValue "G18" "13"
Set-TextValue "D19" "0.3385"
Set-TextValue "E19" "1.29%"
Set-TextValue "G19" "13"
Set-TextValue "D20" "8.212"
Set-TextValue "E20" "-2.50%"
Set-TextValue "G20" "13"
Set-TextValue "D21" "0.1389"
Set-TextValue "E21" "1.04%"
Set-TextValue "G21" "13"
Set-TextValue "E22" "3.68%"
Set-TextValue "G22" "13"
Set-TextValue "D23" "0.001305"
Set-TextValue "E23" "2.99%"
Set-TextValue "G23" "13"
Set-TextValue "D24" "0.004253"
Set-TextValue "E24" "-5.60%"
Set-TextValue "G24" "13"
Set-TextValue "E25" "0.42%"
Set-TextValue "G25" "13"
Set-TextValue "G26" "13"
Set-TextValue "G27" "13"
Set-TextValue "G28" "13"
Set-TextValue "G29" "13"
Set-TextValue "G30" "13"
Set-TextValue "G31" "13"
Set-TextValue "G32" "13"
Set-TextValue "G33" "13"
Set-TextValue "G34" "13"
Set-TextValue "G35" "13"
Set-TextValue "G36" "13"
Set-TextValue "G37" "13"
Set-TextValue "D38" "0.02697"
Set-TextValue "E38" "10.06%"
Set-TextValue "G38" "13"
Set-TextValue "D39" "0.05669"
Set-TextValue "E39" "7.33%"
Set-TextValue "G39" "13"
Set-TextValue "E40" "0.59%"
Set-TextValue "G40" "13"
Set-TextValue "D41" "0.007964"
Set-TextValue "E41" "4.22%"
Set-TextValue "G41" "13"
Set-TextValue "D42" "0.1444"
Set-TextValue "E42" "7.06%"
Set-TextValue "G42" "13"
Set-TextValue "D43" "0.007700"
Set-TextValue "E43" "4.73%"
Set-TextValue "G43" "13"
Set-TextValue "D44" "0.008694"
Set-TextValue "E44" "18.34%"
Set-TextValue "G44" "13"
Set-TextValue "D45" "0.3503"
Set-TextValue "E45" "15.39%"
Set-TextValue "G45" "13"
Set-TextValue "D46" "0.00006875"
Set-TextValue "E46" "4.59%"
Set-TextValue "G46" "13"
Set-TextValue "E47" "0.41%"
Set-TextValue "G47" "13"
Set-TextValue "D48" "0.05546"
Set-TextValue "E48" "-28.17%"
Set-TextValue "G48" "13"
Set-TextValue "E49" "-4.81%"
Set-TextValue "G49" "13"
Set-TextValue "E50" "0.41%"
Set-TextValue "G50" "13"
Set-TextValue "D51" "0.0001999"
Set-TextValue "E51" "0.41%"
Set-TextValue "G51" "13"
